# Journal de Bord - add Q7 Implementation log entries + minor spell-check pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# --- New log rows -------------------------------------------------------
# Row 8 gains a "Q7" marker in the Question-number column (E), matching
# the existing pattern used by rows 2-7 (E2:E7 hold Q1..Q6).
$ws.Range("E8").Value = "Q7"

# Row 33's description is refined from the generic "Implémentation" to the
# more specific "Implémentation Creation de compte".
$ws.Range("C33").Value = "Implémentation Creation de compte"

# New row 34: a continuation line (no time stamp) carrying the old
# "Implémentation" wording, repurposed here as the "Q7" label.
$ws.Range("C34").Value = "Q7"

# New rows 35 and 36: two more implementation entries, each with their own
# timestamp (stored as the usual day-fraction, formatted h:mm).
$ws.Range("B35").Value = 0.40277777777777773
$ws.Range("B35").NumberFormat = "h:mm"
$ws.Range("C35").Value = "Implémentation barre de navigation"

$ws.Range("B36").Value = 0.42708333333333331
$ws.Range("B36").NumberFormat = "h:mm"
$ws.Range("C36").Value = "Implémentation Connexion"

# --- Cosmetic touch-ups ---------------------------------------------------
# Re-fit the columns that received new, longer text so the sheet reads
# cleanly (mirrors the author's manual column-width nudge in the diff).
# Target widths (from the saved file) are ~29.46, ~7.86 and ~7.40
# characters; feed the nearest values this engine's width grid supports.
$ws.Columns("C").ColumnWidth = 28.666666666666668
$ws.Columns("F").ColumnWidth = 7.0
$ws.Columns("G").ColumnWidth = 6.5

# Restore the selection to the last-edited cell, like the saved file shows.
[void]$ws.Range("C33").Select()
